$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-6, i.e. the four member records under the header row)
# are appended again immediately below the existing table (rows 7-11),
# duplicating the member records.
$source = $ws.Range("A2:D6").Value2
$ws.Range("A7:D11").Value2 = $source
